# Auto-update draw results: append the latest Pick 4 draw as a new row
# at the bottom of the Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Find the first empty row right after the existing data (mirrors how the
# upstream exporter appends a new result row each day).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$newRow = $lastRow + 1

$date        = "2025-10-21"
$game        = "Pick 4"
$phase       = "251021"
$result      = "2-7-2-6"
$insertedAt  = "2025-10-21T21:38:44.173+04:00"

# The existing rows store every value (even date-/number-looking ones) as
# plain text, so force Text format on the new row before writing the
# values - this stops Excel from reinterpreting "2025-10-21" / "251021" as
# a date or a number. The style is reset back to Normal afterwards so the
# new cells end up with the same (default) formatting as the rest of the
# sheet.
$newRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5))
$newRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $date
$ws.Cells.Item($newRow, 2).Value = $game
$ws.Cells.Item($newRow, 3).Value = $phase
$ws.Cells.Item($newRow, 4).Value = $result
$ws.Cells.Item($newRow, 5).Value = $insertedAt

$newRange.Style = "Normal"
